$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.015.53'
$ws.Range("E2").Value = '  -3.10%  '
$ws.Range("D3").Value = '3.168.04'
$ws.Range("E3").Value = '  -8.28%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.00%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.604'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '3.165.32'
$ws.Range("E9").Value = '  -8.32%  '
$ws.Range("E10").Value = '  -6.44%  '
$ws.Range("E11").Value = '  -4.84%  '
$ws.Range("E12").Value = '  -4.02%  '
$ws.Range("D13").Value = '3.717.11'
$ws.Range("E13").Value = '  -8.24%  '
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("E15").Value = '  -5.67%  '
$ws.Range("D16").Value = '64.140.02'
$ws.Range("E16").Value = '  -2.97%  '
$ws.Range("E17").Value = '  -5.71%  '
$ws.Range("D18").Value = '3.167.54'
$ws.Range("E18").Value = '  -8.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '351.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.51%  '
$ws.Range("E25").Value = '  -7.00%  '
$ws.Range("E26").Value = '  -4.24%  '
$ws.Range("E27").Value = '  -4.18%  '
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  -3.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.88'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.58'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.19'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.25'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.25%  '
$ws.Range("E37").Value = '  -7.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.802'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -9.10%  '
$ws.Range("E39").Value = '  -10.62%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.36%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.52%  '
$ws.Range("D42").Value = '2.635.39'
$ws.Range("E42").Value = '  -4.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.13'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.56%  '
$ws.Range("E44").Value = '  -7.83%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0651'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '323.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0270'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.10%  '
$ws.Range("E50").Value = '  -0.56%  '
$ws.Range("E51").Value = '  -0.08%  '
